# Insert a new data row at Excel row 37 (pushes the existing row 37..127
# down to 38..128, expanding the used range to A1:R128) and populate it
# with the new "Poroto granado" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(37).Insert()

$ws.Range("A37").Value = 8
$ws.Range("B37").Value = "Terminal La Palmera de La Serena"
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = 45012
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = 100112030
$ws.Range("G37").Value = "Poroto granado"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 500
$ws.Range("K37").Value = 35000
$ws.Range("L37").Value = 36000
$ws.Range("M37").Value = 35500
$ws.Range("N37").Value = "$/malla 25 kilos"
$ws.Range("O37").Value = "Provincia del Elquí"
$ws.Range("P37").Value = 1420
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
